# edit.ps1 -- applies the "Se sube archivo de avanze" changes:
#   1. "...que esta logueado" -> split "logueado" into its own run,
#      wrapped in <w:proofErr w:type="spellStart"/> ... spellEnd (spell-check marker)
#   2. "baja?, eso ocasiona..." -> split "baja?" into its own run,
#      wrapped in <w:proofErr w:type="gramStart"/> ... gramEnd (grammar-check marker)
#   3. Insert a new "DESCUENTOS" heading block (3 paragraphs: blank, title, blank)
#      right before the existing "PRODUCTOS" heading paragraph
#   4. "Al momento de buscar agregar el loader..." -> split "loader" into
#      its own run, wrapped in <w:proofErr w:type="spellStart"/> ... spellEnd

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) "Agregar nombre de usuario y rol del usuario que esta logueado"
# ---------------------------------------------------------------------
$idx1 = Find-ParaIndex $d "logueado"
$p1 = $d.Paragraphs($idx1).Range
$xml1 = "<w:p $wns><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>" +
        "<w:t xml:space=`"preserve`">Agregar nombre de usuario y rol del usuario que esta </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>logueado</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "</w:p>"
$p1.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "baja?, eso ocasiona que cuando le doy guardar se guarde..."
#    The grammar-check marker wraps the whole sentence fragment up to
#    and including "baja?" (gramStart right after <w:pPr>, gramEnd right
#    after the "baja?" run).
# ---------------------------------------------------------------------
$idx2 = Find-ParaIndex $d "eso ocasiona"
$p2 = $d.Paragraphs($idx2).Range
$xml2 = "<w:p $wns><w:pPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>" +
        "<w:proofErr w:type=`"gramStart`"/>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>" +
        "<w:t xml:space=`"preserve`">En apellido materno es correcto que se concatene la palabra </w:t></w:r>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>baja?</w:t></w:r>" +
        "<w:proofErr w:type=`"gramEnd`"/>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>" +
        "<w:t>, eso ocasiona que cuando le doy guardar se guarde la palabra baja en el campo apellido materno</w:t></w:r>" +
        "</w:p>"
$p2.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) Insert "DESCUENTOS" heading block before the "PRODUCTOS" heading
# ---------------------------------------------------------------------
$idx3 = Find-ParaIndex $d "PRODUCTOS"
$p3 = $d.Paragraphs($idx3).Range
$hdrRpr = "<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>"
$xml3 = "<w:p $wns><w:pPr>$hdrRpr</w:pPr></w:p>" +
        "<w:p $wns><w:pPr>$hdrRpr</w:pPr><w:r>$hdrRpr<w:t>DESCUENTOS</w:t></w:r></w:p>" +
        "<w:p $wns><w:pPr>$hdrRpr</w:pPr></w:p>" +
        "<w:p $wns><w:pPr>$hdrRpr</w:pPr><w:r>$hdrRpr<w:t>PRODUCTOS</w:t></w:r></w:p>"
$p3.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) "Al momento de buscar agregar el loader..."
# ---------------------------------------------------------------------
$idx4 = Find-ParaIndex $d "Al momento de buscar"
$p4 = $d.Paragraphs($idx4).Range
$xml4 = "<w:p $wns><w:pPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>" +
        "<w:t xml:space=`"preserve`">Al momento de buscar agregar el </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>loader</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>&#8230;</w:t></w:r>" +
        "</w:p>"
$p4.InsertXML($xml4)

Write-Output "OK"
